$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved (cells store numeric-looking / percent strings as text)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.84%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.722"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-9.58%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05924"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.87%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.663"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8673"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.44%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9376"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.13%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1403"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.51%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.03770"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "8.92%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07107"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.77%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03170"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.33%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09257"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.39%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001543"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.13%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006029"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.50%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006085"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.87%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.498"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.09%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.49%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3145"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.12%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1291"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.17%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.804"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.96%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04226"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.56%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.03%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.39%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.63%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001610"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "34.09%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "1.84%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03823"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.32%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006208"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.75%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1103"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.23%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.30%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.24%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.87%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-11.48%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002438"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "14.33%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
